$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-of-day for Aisyah Dewani Putri's first entry (row 2)
$ws.Range("D2").Value = "19:09:55"

# Row 3 now belongs to a different employee
$ws.Range("A3").Value = "Dinnar Ary Nastiti"
$ws.Range("B3").Value = "MJM005"
$ws.Range("D3").Value = "19:15:25"

# Remove the last record (Fadzli Fiyannuba) entirely - file moved out of module folder
$ws.Rows.Item(4).Delete()
